# Auto-generated Excel COM-interop script applying the market-data refresh
# captured in the commit "chore: update Sheets via scheduled runner".
# For every (sheet, row) below we rewrite the H..N "current price / profit"
# columns to the values recorded after the scheduled run. Plain numeric
# literals only -- the source cells hold static values, not formulas.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 64805.438
$ws.Range("I112").Value = 1849.5
$ws.Range("J112").Value = 73799.14
$ws.Range("K112").Value = 5548.5
$ws.Range("L112").Value = 221397.42
$ws.Range("M112").Value = -4440.5
$ws.Range("N112").Value = -223613.42

$ws.Range("H137").Value = 3761228.8
$ws.Range("I137").Value = 84198.89999999999
$ws.Range("J137").Value = 11115289
$ws.Range("K137").Value = 252596.7
$ws.Range("L137").Value = 33345867
$ws.Range("M137").Value = -250046.7
$ws.Range("N137").Value = -33350967

$ws.Range("H138").Value = 2556.7144
$ws.Range("I138").Value = 826.94116
$ws.Range("J138").Value = 2919.7532
$ws.Range("K138").Value = 2480.82348
$ws.Range("L138").Value = 8759.259600000001
$ws.Range("M138").Value = 2659.17652
$ws.Range("N138").Value = -19039.2596

$ws.Range("H141").Value = 5925.963
$ws.Range("I141").Value = 4761.952
$ws.Range("K141").Value = 14285.856
$ws.Range("M141").Value = -9105.856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4349.278
$ws.Range("I45").Value = 5031.5454
$ws.Range("J45").Value = 4049.08
$ws.Range("K45").Value = 5031.5454
$ws.Range("L45").Value = 4049.08
$ws.Range("M45").Value = -4654.5454
$ws.Range("N45").Value = -4803.08

$ws.Range("H97").Value = 1393.0344
$ws.Range("I97").Value = 1012.9048
$ws.Range("J97").Value = 2390.875
$ws.Range("K97").Value = 1012.9048
$ws.Range("L97").Value = 2390.875
$ws.Range("M97").Value = -516.9048
$ws.Range("N97").Value = -3382.875

$ws.Range("H122").Value = 4498.1133
$ws.Range("I122").Value = 3926.543
$ws.Range("K122").Value = 11779.629
$ws.Range("M122").Value = -9329.629000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1477.8
$ws.Range("I86").Value = 1087.5
$ws.Range("K86").Value = 1087.5
$ws.Range("M86").Value = 35.5

$ws.Range("H89").Value = 1477.8
$ws.Range("I89").Value = 1087.5
$ws.Range("K89").Value = 5437.5
$ws.Range("M89").Value = 178.5

$ws.Range("H94").Value = 295.18182
$ws.Range("I94").Value = 295.18182
$ws.Range("K94").Value = 295.18182
$ws.Range("M94").Value = 155.81818

$ws.Range("H107").Value = 3926.1875
$ws.Range("I107").Value = 3909.3845
$ws.Range("K107").Value = 3909.3845
$ws.Range("M107").Value = -1989.3845

$ws.Range("H134").Value = 4649.8
$ws.Range("I134").Value = 4610.8887
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 13832.6661
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -11297.6661
$ws.Range("N134").Value = -20070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4705.0464
$ws.Range("I31").Value = 2512.2
$ws.Range("J31").Value = 5369.5454
$ws.Range("K31").Value = 2512.2
$ws.Range("L31").Value = 5369.5454
$ws.Range("M31").Value = -2217.2
$ws.Range("N31").Value = -5959.5454

$ws.Range("H34").Value = 4705.0464
$ws.Range("I34").Value = 2512.2
$ws.Range("J34").Value = 5369.5454
$ws.Range("K34").Value = 2512.2
$ws.Range("L34").Value = 5369.5454
$ws.Range("M34").Value = -2310.2
$ws.Range("N34").Value = -5773.5454

$ws.Range("H53").Value = 40000
$ws.Range("J53").Value = 40000
$ws.Range("L53").Value = 40000
$ws.Range("N53").Value = -41214

$ws.Range("H134").Value = 890
$ws.Range("I134").Value = 890
$ws.Range("K134").Value = 2670
$ws.Range("M134").Value = -135

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 209.5
$ws.Range("I40").Value = 80
$ws.Range("K40").Value = 320
$ws.Range("M40").Value = -251

$ws.Range("H92").Value = 1172.9
$ws.Range("J92").Value = 1278.3334
$ws.Range("L92").Value = 3835.0002
$ws.Range("N92").Value = -6331.0002

$ws.Range("H107").Value = 816.65216
$ws.Range("I107").Value = 674.93335
$ws.Range("J107").Value = 1082.375
$ws.Range("K107").Value = 2024.80005
$ws.Range("L107").Value = 3247.125
$ws.Range("M107").Value = -104.8000500000001
$ws.Range("N107").Value = -7087.125

$ws.Range("H131").Value = 1588.92
$ws.Range("I131").Value = 1103.7778
$ws.Range("J131").Value = 1695.4147
$ws.Range("K131").Value = 3311.3334
$ws.Range("L131").Value = 5086.2441
$ws.Range("M131").Value = 1728.6666
$ws.Range("N131").Value = -15166.2441

$ws.Range("H132").Value = 3334466.2
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 2999
$ws.Range("I10").Value = 2999
$ws.Range("K10").Value = 2999
$ws.Range("M10").Value = -2830

$ws.Range("H70").Value = 50947.21
$ws.Range("I70").Value = 122099.8
$ws.Range("K70").Value = 122099.8
$ws.Range("M70").Value = -121829.8

$ws.Range("H73").Value = 50947.21
$ws.Range("I73").Value = 122099.8
$ws.Range("K73").Value = 122099.8
$ws.Range("M73").Value = -121163.8

$ws.Range("H132").Value = 3730.75
$ws.Range("I132").Value = 3521
$ws.Range("K132").Value = 10563
$ws.Range("M132").Value = -8033

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

$ws.Range("H46").Value = 3736.2334
$ws.Range("J46").Value = 4450.75
$ws.Range("L46").Value = 4450.75
$ws.Range("N46").Value = -4826.75

$ws.Range("H58").Value = 5000
$ws.Range("I58").Value = 5000
$ws.Range("K58").Value = 5000
$ws.Range("M58").Value = -4740

$ws.Range("H93").Value = 2526.0557
$ws.Range("J93").Value = 2747.182
$ws.Range("L93").Value = 2747.182
$ws.Range("N93").Value = -5243.182

$ws.Range("H100").Value = 1403.7778
$ws.Range("I100").Value = 1403.7778
$ws.Range("K100").Value = 1403.7778
$ws.Range("M100").Value = -862.7778000000001

$ws.Range("H132").Value = 1416.6154
$ws.Range("I132").Value = 1416.6154
$ws.Range("K132").Value = 4249.8462
$ws.Range("M132").Value = -1719.8462

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 984.0952
$ws.Range("I100").Value = 678
$ws.Range("K100").Value = 1356
$ws.Range("M100").Value = -815

$ws.Range("H104").Value = 91903
$ws.Range("J104").Value = 91903
$ws.Range("L104").Value = 91903
$ws.Range("N104").Value = -98891

$ws.Range("H132").Value = 2296.239
$ws.Range("I132").Value = 1991.3256
$ws.Range("K132").Value = 5973.976799999999
$ws.Range("M132").Value = -3443.976799999999

$ws.Range("H136").Value = 2186.6765
$ws.Range("I136").Value = 1221.5
$ws.Range("K136").Value = 3664.5
$ws.Range("M136").Value = -1114.5

